$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the measured temperature for buffer 2 (F6) from 20 to 22.3
$ws.Range("F6").Value = 22.3

# F5 had a stray/no-op font-applied style (identical in appearance to the
# default style); clear it so the cell goes back to the plain default style.
$ws.Range("F5").ClearFormats()

# Update the selected cell shown in the sheet view
$ws.Range("Q13").Select()
